# Weekly fruit/vegetable price update: insert a new record row right
# before the current row 32, shifting the existing rows 32-46 down to
# rows 33-47 and extending the sheet by one row (A1:R46 -> A1:R47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32; this pushes the old row 32
# (and everything below it) down by one row, preserving formatting.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new data record.
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 45007
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = 100112010
$ws.Cells.Item(32, 7).Value = "Achicoria"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 70
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 7000
$ws.Cells.Item(32, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(32, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(32, 16).Value = 438
$ws.Cells.Item(32, 17).Value = 16
$ws.Cells.Item(32, 18).Value = "Hortaliza"
